# Hometask.docx update:
#  - "GradientBoosting_tuning_" -> "XGBoost_tuning_" (both occurrences)
#  - the now-stale "_GoBack" bookmark (left over from the previous save,
#    sitting in the empty paragraph right under the title) is replaced by
#    a fresh one at the spot where the author's cursor ended up after the
#    last edit, in the middle of "решите задачу"
#
# Word re-splits/re-merges runs as a side effect of these edits, which the
# canonical OOXML reflects; we reproduce that by doing the same edits via
# Find/Replace (which coalesces touched, identically-formatted runs) and
# by dropping bookmarks at the split points (which forces a run break).

$d = $word.ActiveDocument

# 1) Rename the notebook everywhere it is mentioned.
$d.Content.Find.Execute("GradientBoosting_tuning_", $true, $false, $false, $false, $false, $true, 1, $false, "XGBoost_tuning_", 2)

# 2) A couple of incidental run merges that happen elsewhere in the same
#    paragraph once Word resaves the document.
$d.Content.Find.Execute("2) ", $true, $false, $false, $false, $false, $true, 1, $false, "2) ", 2)
$d.Content.Find.Execute("test_medium.csv (", $true, $false, $false, $false, $false, $true, 1, $false, "test_medium.csv (", 2)

# 3) Drop the "_GoBack" bookmark where the author's cursor last was, right
#    in the middle of "решите задачу" ("решите з" | "адачу"). Adding a new
#    "_GoBack" bookmark automatically removes the old (single, special)
#    one that Word maintains, and also forces a run break at this point.
$r = $d.Content
$r.Find.Execute("Затем решите з", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# 4) Merge the runs around the new bookmark (and a later, unrelated, run
#    pair) back together the way Word leaves them after saving.
$d.Content.Find.Execute(". Затем решите з", $true, $false, $false, $false, $false, $true, 1, $false, ". Затем решите з", 2)
$d.Content.Find.Execute("адачу с помощью ", $true, $false, $false, $false, $false, $true, 1, $false, "адачу с помощью ", 2)
$d.Content.Find.Execute(". Сравните три алгоритма ", $true, $false, $false, $false, $false, $true, 1, $false, ". Сравните три алгоритма ", 2)

# 5) "XGBoost_tuning_" (both occurrences) ends up split into two runs,
#    "XGBoost" and "_tuning_" -- only the first word was actually retyped,
#    the "_tuning_" suffix is the untouched remainder of the old run.
#    Force that split by dropping a temporary bookmark at the boundary
#    (which breaks the run) and removing the bookmark again.
$scan = $d.Content
$more = $true
while ($more) {
    $more = $scan.Find.Execute("XGBoost_tuning_", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($more) {
        $splitAt = $scan.Start + 7
        $splitRange = $d.Range($splitAt, $splitAt)
        $d.Bookmarks.Add("TEMP_SPLIT_MARKER", $splitRange)
        $d.Bookmarks("TEMP_SPLIT_MARKER").Delete()
        $scan.Start = $scan.End
        $scan.End = $d.Content.End
    }
}
